$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
  0.11035056273264843,
  -0.0059999999653186364,
  -0.0039999999719491086,
  -0.0079999999474598127,
  -0.0029999999741292527,
  -0.0019999999753483877,
  -0.0099999999293198805,
  -0.0099999999278321816,
  -0.0019999999719728656,
  0.019519257162420089,
  -0.0029999999650511811,
  -0.0034999999613560817,
  -0.0034999999576763585,
  -0.007999999931158186,
  -0.00099999997025967957,
  -0.001999999963985033,
  -0.0019999999630835319,
  -0.003999999951511235,
  -0.0039999999765703009,
  -0.0039999999763935534,
  -0.0039999999765552019,
  0.010013515737603562,
  -0.0049999999642738047,
  -0.019999999875334851,
  -0.019999999873585139,
  -0.0024999999723611666,
  -0.0024999999722465915,
  -0.0019999999747781771,
  -0.0062772939788375837,
  0.028976408432102474,
  -0.0069999999444565475,
  -0.0099999999270234952,
  -0.0039999999609374726
)

for ($i = 0; $i -lt $values.Length; $i++) {
  $row = $i + 1
  $ws.Cells.Item($row, 1).Value = $values[$i]
}
